$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 244; existing rows 244-254 shift down to 246-256.
$ws.Rows("244:245").Insert()

# --- New row 244 ---
$ws.Range("A244").Value2 = 3
$ws.Range("B244").Value2 = "Femacal de La Calera"
$ws.Range("C244").Value2 = "Coquimbo"
$ws.Range("D244").Value2 = 44516
$ws.Range("E244").Value2 = 5
$ws.Range("F244").Value2 = 100112013
$ws.Range("G244").Value2 = "Alcachofa"
$ws.Range("H244").Value2 = "Española"
$ws.Range("I244").Value2 = "Primera"
$ws.Range("J244").Value2 = 12500
$ws.Range("K244").Value2 = 250
$ws.Range("L244").Value2 = 270
$ws.Range("M244").Value2 = 260
$ws.Range("N244").Value2 = "$/unidad"
$ws.Range("O244").Value2 = "Provincia de Quillota"
$ws.Range("P244").Value2 = 260
$ws.Range("Q244").Value2 = 1
$ws.Range("R244").Value2 = "Hortaliza"

# --- New row 245 ---
$ws.Range("A245").Value2 = 3
$ws.Range("B245").Value2 = "Femacal de La Calera"
$ws.Range("C245").Value2 = "Coquimbo"
$ws.Range("D245").Value2 = 44516
$ws.Range("E245").Value2 = 5
$ws.Range("F245").Value2 = 100112013
$ws.Range("G245").Value2 = "Alcachofa"
$ws.Range("H245").Value2 = "Española"
$ws.Range("I245").Value2 = "Segunda"
$ws.Range("J245").Value2 = 5800
$ws.Range("K245").Value2 = 180
$ws.Range("L245").Value2 = 180
$ws.Range("M245").Value2 = 180
$ws.Range("N245").Value2 = "$/unidad"
$ws.Range("O245").Value2 = "Provincia de Quillota"
$ws.Range("P245").Value2 = 180
$ws.Range("Q245").Value2 = 1
$ws.Range("R245").Value2 = "Hortaliza"
